$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns remain text so values like "1.024" or "0.5150" are preserved exactly
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "28.527.02"
$ws.Range("E2").Value = "  +0.89%  "

$ws.Range("D3").Value = "1.876.68"
$ws.Range("E3").Value = "  -0.42%  "

$ws.Range("D4").Value = "1.024"
$ws.Range("E4").Value = "  +1.67%  "

$ws.Range("D5").Value = "317.61"
$ws.Range("E5").Value = "  +1.02%  "

$ws.Range("D6").Value = "1.023"
$ws.Range("E6").Value = "  +1.66%  "

$ws.Range("D7").Value = "0.5150"
$ws.Range("E7").Value = "  +0.22%  "

$ws.Range("D8").Value = "0.3936"
$ws.Range("E8").Value = "  +0.39%  "

$ws.Range("E9").Value = "  -0.25%  "

$ws.Range("E10").Value = "  -0.50%  "

$ws.Range("E11").Value = "  +0.37%  "

$ws.Range("D12").Value = "1.871.51"
$ws.Range("E12").Value = "  -0.79%  "

$ws.Range("E13").Value = "  -1.63%  "

$ws.Range("B14").Value = "BinanceUSD"
$ws.Range("C14").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D14").Value = "1.027"
$ws.Range("E14").Value = "  +2.00%  "

$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").Value = "7.260"
$ws.Range("E15").Value = "  -0.45%  "

$ws.Range("E16").Value = "  +0.30%  "

$ws.Range("D17").Value = "91.54"
$ws.Range("E17").Value = "  -0.03%  "

$ws.Range("D18").Value = "0.06761"
$ws.Range("E18").Value = "  +1.37%  "

$ws.Range("D19").Value = "1.023"
$ws.Range("E19").Value = "  +1.73%  "

$ws.Range("D20").Value = "17.67"
$ws.Range("E20").Value = "  -0.97%  "

$ws.Range("D21").Value = "5.971"
$ws.Range("E21").Value = "  -1.36%  "

$ws.Range("D22").Value = "28.541.13"
$ws.Range("E22").Value = "  +0.73%  "

$ws.Range("D23").Value = "11.19"
$ws.Range("E23").Value = "  +0.08%  "

$ws.Range("D24").Value = "2.273"
$ws.Range("E24").Value = "  -0.30%  "

$ws.Range("D25").Value = "2.086.55"
$ws.Range("E25").Value = "  -0.47%  "

$ws.Range("D26").Value = "161.97"
$ws.Range("E26").Value = "  +1.73%  "

$ws.Range("D27").Value = "20.72"
$ws.Range("E27").Value = "  +0.14%  "

$ws.Range("D28").Value = "2.370"
$ws.Range("E28").Value = "  -6.25%  "

$ws.Range("D29").Value = "127.41"
$ws.Range("E29").Value = "  +1.28%  "

$ws.Range("D30").Value = "0.1052"
$ws.Range("E30").Value = "  -1.48%  "

$ws.Range("D31").Value = "1.036"
$ws.Range("E31").Value = "  -1.51%  "

$ws.Range("D32").Value = "5.856"
$ws.Range("E32").Value = "  -0.77%  "

$ws.Range("D33").Value = "3.668"
$ws.Range("E33").Value = "  +1.70%  "

$ws.Range("D34").Value = "0.02443"
$ws.Range("E34").Value = "  -1.07%  "

$ws.Range("D35").Value = "0.06494"
$ws.Range("E35").Value = "  -1.60%  "

$ws.Range("D36").Value = "9.152"
$ws.Range("E36").Value = "  -6.28%  "

$ws.Range("D37").Value = "0.2184"
$ws.Range("E37").Value = "  -0.64%  "

$ws.Range("D38").Value = "1.250"
$ws.Range("E38").Value = "  +1.56%  "

$ws.Range("B39").Value = "TheSandbox"
$ws.Range("C39").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D39").Value = "0.6455"
$ws.Range("E39").Value = "  -1.35%  "

$ws.Range("B40").Value = "ARBITRUM"
$ws.Range("C40").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D40").Value = "1.188"
$ws.Range("E40").Value = "  -2.18%  "

$ws.Range("D41").Value = "4.983"
$ws.Range("E41").Value = "  -0.87%  "

$ws.Range("D42").Value = "11.18"
$ws.Range("E42").Value = "  -1.23%  "

$ws.Range("D43").Value = "0.6031"
$ws.Range("E43").Value = "  -2.28%  "

$ws.Range("B44").Value = "PancakeSwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D44").Value = "3.722"
$ws.Range("E44").Value = "  +0.95%  "

$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "12.92"
$ws.Range("E45").Value = "  -2.21%  "

$ws.Range("D46").Value = "1.278"
$ws.Range("E46").Value = "  -0.77%  "

$ws.Range("E47").Value = "  -1.39%  "

$ws.Range("D48").Value = "1.215"
$ws.Range("E48").Value = "  -1.95%  "

$ws.Range("D49").Value = "122.08"
$ws.Range("E49").Value = "  +0.29%  "

$ws.Range("D50").Value = "0.06869"
$ws.Range("E50").Value = "  -0.44%  "

$ws.Range("D51").Value = "76.32"
$ws.Range("E51").Value = "  -3.64%  "
